# Update the "want to go" counts (F column) for two events that exist on
# both the "展览" sheet and the "全部类型" sheet.
#   Row 3 (丽水·龙泉ACG动漫游戏博览会): F3 1300 -> 1303
#   Row 8 (丽水·AEO纯白礼赞动漫嘉年华):   F8 171  -> 179

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1303
    $ws.Range("F8").Value = 179
}
